$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the marked scores for the two rectangular contours tracked
$ws.Range("B2").Value = 0
$ws.Range("B6").Value = 2

# Move the active selection to B6, matching the final cursor position
$ws.Range("B6").Select()
